$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.117659
$ws.Range("H2").Value = 12.352977
$ws.Range("I2").Value = 0.411783455701328
$ws.Range("J2").Value = 0.4117834557013281
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 70.68591675326999
$ws.Range("R2").Value = 636.17325077943
$ws.Range("S2").Value = 0.02307409640158157
$ws.Range("T2").Value = 0.02307409640158158

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.117659
$ws.Range("H3").Value = 12.352977
$ws.Range("I3").Value = 0.411783455701328
$ws.Range("J3").Value = 0.4117834557013281
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 1055.950201854548
$ws.Range("R3").Value = 9503.551816690931
$ws.Range("S3").Value = 0.3446952076452232
$ws.Range("T3").Value = 0.3446952076452232

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.117659
$ws.Range("H4").Value = 12.352977
$ws.Range("I4").Value = 0.411783455701328
$ws.Range("J4").Value = 0.4117834557013281
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 134.834344351798
$ws.Range("R4").Value = 1213.509099166182
$ws.Range("S4").Value = 0.04401415165452327
$ws.Range("T4").Value = 0.04401415165452328

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.004148
$ws.Range("H5").Value = 9.012444
$ws.Range("I5").Value = 0.3004276082303642
$ws.Range("J5").Value = 0.3004276082303643
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 51.57079676644
$ws.Range("R5").Value = 464.13717089796
$ws.Range("S5").Value = 0.01683432274421424
$ws.Range("T5").Value = 0.01683432274421425

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.004148
$ws.Range("H6").Value = 9.012444
$ws.Range("I6").Value = 0.3004276082303642
$ws.Range("J6").Value = 0.3004276082303643
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 770.3966469785227
$ws.Range("R6").Value = 6933.569822806704
$ws.Range("S6").Value = 0.2514815866629515
$ws.Range("T6").Value = 0.2514815866629515

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.004148
$ws.Range("H7").Value = 9.012444
$ws.Range("I7").Value = 0.3004276082303642
$ws.Range("J7").Value = 0.3004276082303643
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 98.37199387218934
$ws.Range("R7").Value = 885.347944849704
$ws.Range("S7").Value = 0.03211169882319852
$ws.Range("T7").Value = 0.03211169882319852

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.877766666666667
$ws.Range("H8").Value = 8.6333
$ws.Range("I8").Value = 0.2877889360683077
$ws.Range("J8").Value = 0.2877889360683077
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 49.40126781633332
$ws.Range("R8").Value = 444.611410347
$ws.Range("S8").Value = 0.01612612056703207
$ws.Range("T8").Value = 0.01612612056703208

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.877766666666667
$ws.Range("H9").Value = 8.6333
$ws.Range("I9").Value = 0.2877889360683077
$ws.Range("J9").Value = 0.2877889360683077
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 737.9868737447555
$ws.Range("R9").Value = 6641.8818637028
$ws.Range("S9").Value = 0.2409020219307058
$ws.Range("T9").Value = 0.2409020219307059

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.877766666666667
$ws.Range("H10").Value = 8.6333
$ws.Range("I10").Value = 0.2877889360683077
$ws.Range("J10").Value = 0.2877889360683077
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 94.23358799197777
$ws.Range("R10").Value = 848.1022919278
$ws.Range("S10").Value = 0.03076079357056973
$ws.Range("T10").Value = 0.03076079357056974

